$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 10 (Auto Connect module / General settings) test-case summary figures
$ws.Range("C62").Value = 738
$ws.Range("C63").Value = 926
$ws.Range("C64").Value = 636

# Move the active selection/scroll position down to the newly-filled rows
$ws.Range("C64").Select() | Out-Null
